$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.011274565685031
$ws.Cells.Item(2, 4).Value = 1.02853621527831
$ws.Cells.Item(2, 5).Value = 1.013456855504116
$ws.Cells.Item(2, 6).Value = 1.024805272200165
$ws.Cells.Item(2, 9).Value = 1.03066174483276
$ws.Cells.Item(2, 10).Value = 1.016523982702583
$ws.Cells.Item(2, 11).Value = 1.031352746829687
$ws.Cells.Item(2, 12).Value = 1.0163179513266
$ws.Cells.Item(2, 13).Value = 1.027632697022692
$ws.Cells.Item(2, 14).Value = 1.017967563354951
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.012502400820087
$ws.Cells.Item(3, 4).Value = 1.028995235621454
$ws.Cells.Item(3, 5).Value = 1.014505626989326
$ws.Cells.Item(3, 6).Value = 1.02625049380944
$ws.Cells.Item(3, 9).Value = 1.03077321492864
$ws.Cells.Item(3, 10).Value = 1.017383520173195
$ws.Cells.Item(3, 11).Value = 1.031620664893503
$ws.Cells.Item(3, 12).Value = 1.017170776662011
$ws.Cells.Item(3, 13).Value = 1.028883354028757
$ws.Cells.Item(3, 14).Value = 1.018828321467362
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.013296300454315
$ws.Cells.Item(4, 4).Value = 1.029289502745805
$ws.Cells.Item(4, 5).Value = 1.015184093410783
$ws.Cells.Item(4, 6).Value = 1.027181442250794
$ws.Cells.Item(4, 9).Value = 1.030841896708439
$ws.Cells.Item(4, 10).Value = 1.017938678159727
$ws.Cells.Item(4, 11).Value = 1.031790556900791
$ws.Cells.Item(4, 12).Value = 1.017721864839092
$ws.Cells.Item(4, 13).Value = 1.029687914664434
$ws.Cells.Item(4, 14).Value = 1.019384267841912
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.013629917383936
$ws.Cells.Item(5, 4).Value = 1.029412552283105
$ws.Cells.Item(5, 5).Value = 1.015469284551913
$ws.Cells.Item(5, 6).Value = 1.027571809621712
$ws.Cells.Item(5, 9).Value = 1.030869944160312
$ws.Cells.Item(5, 10).Value = 1.018171824287327
$ws.Cells.Item(5, 11).Value = 1.031861147344941
$ws.Cells.Item(5, 12).Value = 1.017953365290317
$ws.Cells.Item(5, 13).Value = 1.030025029930159
$ws.Cells.Item(5, 14).Value = 1.019617745063755
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.013685925106522
$ws.Cells.Item(6, 4).Value = 1.029433174056169
$ws.Cells.Item(6, 5).Value = 1.015517167334347
$ws.Cells.Item(6, 6).Value = 1.027637295289126
$ws.Cells.Item(6, 9).Value = 1.030874604956507
$ws.Cells.Item(6, 10).Value = 1.018210956396288
$ws.Cells.Item(6, 11).Value = 1.031872950942328
$ws.Cells.Item(6, 12).Value = 1.017992224905394
$ws.Cells.Item(6, 13).Value = 1.030081567346615
$ws.Cells.Item(6, 14).Value = 1.0196569327448
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.013300758800292
$ws.Cells.Item(7, 4).Value = 1.029291149537283
$ws.Cells.Item(7, 5).Value = 1.015187904288528
$ws.Cells.Item(7, 6).Value = 1.027186662297118
$ws.Cells.Item(7, 9).Value = 1.030842274728561
$ws.Cells.Item(7, 10).Value = 1.017941794419407
$ws.Cells.Item(7, 11).Value = 1.031791503405852
$ws.Cells.Item(7, 12).Value = 1.017724958852899
$ws.Cells.Item(7, 13).Value = 1.029692423621076
$ws.Cells.Item(7, 14).Value = 1.019387388527038
$ws.Cells.Item(8, 2).Value = 1.019999999999999
$ws.Cells.Item(8, 3).Value = 1.011689642790092
$ws.Cells.Item(8, 4).Value = 1.028691911481741
$ws.Cells.Item(8, 5).Value = 1.013811326775708
$ws.Cells.Item(8, 6).Value = 1.025294564878773
$ws.Cells.Item(8, 9).Value = 1.030700129662548
$ws.Cells.Item(8, 10).Value = 1.0168146798471
$ws.Cells.Item(8, 11).Value = 1.031444008255411
$ws.Cells.Item(8, 12).Value = 1.016606323357727
$ws.Cells.Item(8, 13).Value = 1.028056336844365
$ws.Cells.Item(8, 14).Value = 1.018258673322757
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.008845973171714
$ws.Cells.Item(9, 4).Value = 1.027615005125369
$ws.Cells.Item(9, 5).Value = 1.011384318510748
$ws.Cells.Item(9, 6).Value = 1.021928057433384
$ws.Cells.Item(9, 9).Value = 1.030423302941615
$ws.Cells.Item(9, 10).Value = 1.014820657002179
$ws.Cells.Item(9, 11).Value = 1.030805180663596
$ws.Cells.Item(9, 12).Value = 1.014629343716847
$ws.Cells.Item(9, 13).Value = 1.025137233537864
$ws.Cells.Item(9, 14).Value = 1.016261818736679
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.006946842443587
$ws.Cells.Item(10, 4).Value = 1.0268830783563
$ws.Cells.Item(10, 5).Value = 1.009765312790407
$ws.Cells.Item(10, 6).Value = 1.019661688225166
$ws.Cells.Item(10, 9).Value = 1.030221098410078
$ws.Cells.Item(10, 10).Value = 1.013485866536836
$ws.Cells.Item(10, 11).Value = 1.030361570022238
$ws.Cells.Item(10, 12).Value = 1.013307340221098
$ws.Cells.Item(10, 13).Value = 1.023166695640034
$ws.Cells.Item(10, 14).Value = 1.014925132715772
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.006123656106344
$ws.Cells.Item(11, 4).Value = 1.026562849577683
$ws.Cells.Item(11, 5).Value = 1.009063999406542
$ws.Cells.Item(11, 6).Value = 1.018675039880131
$ws.Cells.Item(11, 9).Value = 1.030129366162488
$ws.Cells.Item(11, 10).Value = 1.012906567117103
$ws.Cells.Item(11, 11).Value = 1.030165294451493
$ws.Cells.Item(11, 12).Value = 1.012733918284392
$ws.Cells.Item(11, 13).Value = 1.022307581042747
$ws.Cells.Item(11, 14).Value = 1.014345010624417
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.005817756154902
$ws.Cells.Item(12, 4).Value = 1.026443408192872
$ws.Cells.Item(12, 5).Value = 1.008803456746802
$ws.Cells.Item(12, 6).Value = 1.018307753649629
$ws.Cells.Item(12, 9).Value = 1.03009466606873
$ws.Cells.Item(12, 10).Value = 1.012691187580854
$ws.Cells.Item(12, 11).Value = 1.030091760776301
$ws.Cells.Item(12, 12).Value = 1.012520773318444
$ws.Cells.Item(12, 13).Value = 1.021987583092373
$ws.Cells.Item(12, 14).Value = 1.014129325224522
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.005883378768562
$ws.Cells.Item(13, 4).Value = 1.026469051112548
$ws.Cells.Item(13, 5).Value = 1.008859346048697
$ws.Cells.Item(13, 6).Value = 1.018386574134373
$ws.Cells.Item(13, 9).Value = 1.030102137701758
$ws.Cells.Item(13, 10).Value = 1.012737396403895
$ws.Cells.Item(13, 11).Value = 1.030107562407695
$ws.Cells.Item(13, 12).Value = 1.012566500488306
$ws.Cells.Item(13, 13).Value = 1.022056263837041
$ws.Cells.Item(13, 14).Value = 1.014175599669391
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.006098373035732
$ws.Cells.Item(14, 4).Value = 1.026552986587392
$ws.Cells.Item(14, 5).Value = 1.009042463774448
$ws.Cells.Item(14, 6).Value = 1.018644696260421
$ws.Cells.Item(14, 9).Value = 1.030126510621029
$ws.Cells.Item(14, 10).Value = 1.012888767912514
$ws.Cells.Item(14, 11).Value = 1.030159228943182
$ws.Cells.Item(14, 12).Value = 1.012716302739939
$ws.Cells.Item(14, 13).Value = 1.022281147982192
$ws.Cells.Item(14, 14).Value = 1.014327186142916
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.006230820497123
$ws.Cells.Item(15, 4).Value = 1.026604636568748
$ws.Cells.Item(15, 5).Value = 1.009155282787606
$ws.Cells.Item(15, 6).Value = 1.018803627499988
$ws.Cells.Item(15, 9).Value = 1.03014144456633
$ws.Cells.Item(15, 10).Value = 1.012982006055535
$ws.Cells.Item(15, 11).Value = 1.030190979194253
$ws.Cells.Item(15, 12).Value = 1.012808580831506
$ws.Cells.Item(15, 13).Value = 1.022419589184444
$ws.Cells.Item(15, 14).Value = 1.014420556694794
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.007001456342757
$ws.Cells.Item(16, 4).Value = 1.026904261502182
$ws.Cells.Item(16, 5).Value = 1.009811850672825
$ws.Cells.Item(16, 6).Value = 1.019727056788365
$ws.Cells.Item(16, 9).Value = 1.030227098455601
$ws.Cells.Item(16, 10).Value = 1.013524284559835
$ws.Cells.Item(16, 11).Value = 1.030374508026889
$ws.Cells.Item(16, 12).Value = 1.013345375371713
$ws.Cells.Item(16, 13).Value = 1.023223588443539
$ws.Cells.Item(16, 14).Value = 1.014963605296771
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.007484625301137
$ws.Cells.Item(17, 4).Value = 1.027091326033321
$ws.Cells.Item(17, 5).Value = 1.010223623390185
$ws.Cells.Item(17, 6).Value = 1.020304877972666
$ws.Cells.Item(17, 9).Value = 1.030279709329392
$ws.Cells.Item(17, 10).Value = 1.013864084565271
$ws.Cells.Item(17, 11).Value = 1.030488510105351
$ws.Cells.Item(17, 12).Value = 1.01368182679472
$ws.Cells.Item(17, 13).Value = 1.023726344380258
$ws.Cells.Item(17, 14).Value = 1.01530388785719
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.007766367832308
$ws.Cells.Item(18, 4).Value = 1.027200119140072
$ws.Cells.Item(18, 5).Value = 1.010463777265698
$ws.Cells.Item(18, 6).Value = 1.020641400734278
$ws.Cells.Item(18, 9).Value = 1.030309993355115
$ws.Cells.Item(18, 10).Value = 1.014062156370213
$ws.Cells.Item(18, 11).Value = 1.030554601313496
$ws.Cells.Item(18, 12).Value = 1.013877978279084
$ws.Cells.Item(18, 13).Value = 1.024019028436331
$ws.Cells.Item(18, 14).Value = 1.015502240946814
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.007862421013543
$ws.Cells.Item(19, 4).Value = 1.027237160710456
$ws.Cells.Item(19, 5).Value = 1.010545659167395
$ws.Cells.Item(19, 6).Value = 1.020756059863573
$ws.Cells.Item(19, 9).Value = 1.030320251042939
$ws.Cells.Item(19, 10).Value = 1.014129672141652
$ws.Cells.Item(19, 11).Value = 1.030577068097479
$ws.Cells.Item(19, 12).Value = 1.013944844835894
$ws.Cells.Item(19, 13).Value = 1.024118730434403
$ws.Cells.Item(19, 14).Value = 1.015569852598392
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.007432794331597
$ws.Cells.Item(20, 4).Value = 1.027071288709237
$ws.Cells.Item(20, 5).Value = 1.010179446802759
$ws.Cells.Item(20, 6).Value = 1.020242936079603
$ws.Cells.Item(20, 9).Value = 1.030274106359939
$ws.Cells.Item(20, 10).Value = 1.013827640478709
$ws.Cells.Item(20, 11).Value = 1.030476320569007
$ws.Cells.Item(20, 12).Value = 1.013645738579692
$ws.Cells.Item(20, 13).Value = 1.023672461897711
$ws.Cells.Item(20, 14).Value = 1.015267392015845
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.006035066257158
$ws.Cells.Item(21, 4).Value = 1.026528283310985
$ws.Cells.Item(21, 5).Value = 1.008988541394068
$ws.Cells.Item(21, 6).Value = 1.018568707887153
$ws.Cells.Item(21, 9).Value = 1.030119350696332
$ws.Cells.Item(21, 10).Value = 1.012844198372178
$ws.Cells.Item(21, 11).Value = 1.030144031768088
$ws.Cells.Item(21, 12).Value = 1.01267219388557
$ws.Cells.Item(21, 13).Value = 1.022214949659526
$ws.Cells.Item(21, 14).Value = 1.014282553308719
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.005155493957107
$ws.Cells.Item(22, 4).Value = 1.026184015943033
$ws.Cells.Item(22, 5).Value = 1.008239518082255
$ws.Cells.Item(22, 6).Value = 1.017511416401771
$ws.Cells.Item(22, 9).Value = 1.030018424404569
$ws.Cells.Item(22, 10).Value = 1.012224699584895
$ws.Cells.Item(22, 11).Value = 1.029931474640954
$ws.Cells.Item(22, 12).Value = 1.012059215563109
$ws.Cells.Item(22, 13).Value = 1.021293433969382
$ws.Cells.Item(22, 14).Value = 1.013662174762102
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.005621845804848
$ws.Cells.Item(23, 4).Value = 1.026366789017498
$ws.Cells.Item(23, 5).Value = 1.008636614445017
$ws.Cells.Item(23, 6).Value = 1.018072348104115
$ws.Cells.Item(23, 9).Value = 1.030072270722477
$ws.Cells.Item(23, 10).Value = 1.012553219388124
$ws.Cells.Item(23, 11).Value = 1.030044499273701
$ws.Cells.Item(23, 12).Value = 1.012384250459518
$ws.Cells.Item(23, 13).Value = 1.021782433659608
$ws.Cells.Item(23, 14).Value = 1.013991161101133
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.007456214768037
$ws.Cells.Item(24, 4).Value = 1.027080343697817
$ws.Cells.Item(24, 5).Value = 1.010199408381864
$ws.Cells.Item(24, 6).Value = 1.020270926531837
$ws.Cells.Item(24, 9).Value = 1.030276639346022
$ws.Cells.Item(24, 10).Value = 1.013844108387407
$ws.Cells.Item(24, 11).Value = 1.030481829745161
$ws.Cells.Item(24, 12).Value = 1.013662045583888
$ws.Cells.Item(24, 13).Value = 1.023696810817492
$ws.Cells.Item(24, 14).Value = 1.015283883310861
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.009581705315051
$ws.Cells.Item(25, 4).Value = 1.027895884775391
$ws.Cells.Item(25, 5).Value = 1.012011925387101
$ws.Cells.Item(25, 6).Value = 1.02280224431607
$ws.Cells.Item(25, 9).Value = 1.030497985149791
$ws.Cells.Item(25, 10).Value = 1.015337109338483
$ws.Cells.Item(25, 11).Value = 1.030973463410629
$ws.Cells.Item(25, 12).Value = 1.01514113970018
$ws.Cells.Item(25, 13).Value = 1.025896190058279
$ws.Cells.Item(25, 14).Value = 1.016779004494538
